# Update "想去人数" (want-to-go count) values in both the "展览" sheet
# and the mirrored rows in the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - rows 2-9, column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 5370
$wsExpo.Range("F3").Value = 584
$wsExpo.Range("F4").Value = 11317
$wsExpo.Range("F5").Value = 278
$wsExpo.Range("F6").Value = 586
$wsExpo.Range("F7").Value = 163
$wsExpo.Range("F8").Value = 247
$wsExpo.Range("F9").Value = 972

# Sheet "全部类型" (All types) - the same events appear at different rows
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 5370
$wsAll.Range("F5").Value = 584
$wsAll.Range("F7").Value = 11317
$wsAll.Range("F8").Value = 278
$wsAll.Range("F9").Value = 586
$wsAll.Range("F10").Value = 163
$wsAll.Range("F13").Value = 247
$wsAll.Range("F14").Value = 972
